{"js": "// The document is a single 20x5 table of arithmetic equations\n// (\"a+b=c\" / \"a-b=c\") \u2014 one equation per cell, one run/paragraph per\n// cell. The commit replaces each equation's text with a new one,\n// strictly in document (row-major) order. Some \"old\" equations repeat\n// (e.g. \"63-47=16\" appears twice), so we must replace by position,\n// not by searching for old text.\nconst newValues = [\"80-16=64\", \"63-49=14\", \"28+43=71\", \"47+47=94\", \"83-17=66\", \"59+12=71\", \"35+8=43\", \"3+89=92\", \"25+19=44\", \"3+38=41\", \"53-28=25\", \"92-23=69\", \"76+5=81\", \"81-2=79\", \"45-37=8\", \"8+18=26\", \"63-9=54\", \"68+6=74\", \"52-9=43\", \"94-46=48\", \"18+73=91\", \"25-18=7\", \"79+12=91\", \"94-26=68\", \"60-16=44\", \"49+6=55\", \"97-19=78\", \"85+7=92\", \"9+29=38\", \"49+8=57\", \"89+4=93\", \"15+17=32\", \"49+34=83\", \"18+9=27\", \"42+19=61\", \"30-24=6\", \"8+54=62\", \"44+19=63\", \"77+16=93\", \"48+5=53\", \"68-39=29\", \"17+14=31\", \"51-22=29\", \"96-9=87\", \"71-4=67\", \"29+67=96\", \"15+78=93\", \"8+4=12\", \"75-69=6\", \"22-14=8\", \"34-6=28\", \"64-39=25\", \"32+9=41\", \"15+46=61\", \"90-62=28\", \"58+9=67\", \"62-45=17\", \"14+8=22\", \"58+39=97\", \"52-38=14\", \"53-34=19\", \"35-18=17\", \"84-46=38\", \"73-6=67\", \"68-19=49\", \"70-51=19\", \"24+39=63\", \"60-53=7\", \"26+69=95\", \"49+25=74\", \"90-44=46\", \"92-23=69\", \"9+75=84\", \"9+74=83\", \"35-26=9\", \"90-35=55\", \"84-59=25\", \"94-16=78\", \"47-9=38\", \"28+68=96\", \"78+9=87\", \"76-29=47\", \"51-18=33\", \"27+46=73\", \"38+37=75\", \"57+36=93\", \"2+89=91\", \"57+7=64\", \"93-39=54\", \"60-5=55\", \"66+25=91\", \"43-4=39\", \"26+59=85\", \"25+17=42\", \"90-41=49\", \"47-18=29\", \"53-49=4\", \"37+57=94\", \"58+24=82\", \"26+69=95\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// table.values is row-major: [[r0c0, r0c1, ...], [r1c0, ...], ...]\nconst rows = table.values;\nlet idx = 0;\nfor (let r = 0; r < rows.length; r++) {\n  for (let c = 0; c < rows[r].length; c++) {\n    rows[r][c] = newValues[idx];\n    idx++;\n  }\n}\n\ntable.values = rows;\nawait context.sync();\n", "ps1": "# The document is a single 20x5 table of arithmetic equations\n# (\"a+b=c\" / \"a-b=c\") \u2014 one equation per cell, one run/paragraph per\n# cell. The commit replaces each equation's text with a new one,\n# strictly in document (row-major) order. Some \"old\" equations repeat\n# (e.g. \"63-47=16\" appears twice), so we must replace by position,\n# not by searching for old text.\n$newValues = @(\"80-16=64\",\"63-49=14\",\"28+43=71\",\"47+47=94\",\"83-17=66\",\"59+12=71\",\"35+8=43\",\"3+89=92\",\"25+19=44\",\"3+38=41\",\"53-28=25\",\"92-23=69\",\"76+5=81\",\"81-2=79\",\"45-37=8\",\"8+18=26\",\"63-9=54\",\"68+6=74\",\"52-9=43\",\"94-46=48\",\"18+73=91\",\"25-18=7\",\"79+12=91\",\"94-26=68\",\"60-16=44\",\"49+6=55\",\"97-19=78\",\"85+7=92\",\"9+29=38\",\"49+8=57\",\"89+4=93\",\"15+17=32\",\"49+34=83\",\"18+9=27\",\"42+19=61\",\"30-24=6\",\"8+54=62\",\"44+19=63\",\"77+16=93\",\"48+5=53\",\"68-39=29\",\"17+14=31\",\"51-22=29\",\"96-9=87\",\"71-4=67\",\"29+67=96\",\"15+78=93\",\"8+4=12\",\"75-69=6\",\"22-14=8\",\"34-6=28\",\"64-39=25\",\"32+9=41\",\"15+46=61\",\"90-62=28\",\"58+9=67\",\"62-45=17\",\"14+8=22\",\"58+39=97\",\"52-38=14\",\"53-34=19\",\"35-18=17\",\"84-46=38\",\"73-6=67\",\"68-19=49\",\"70-51=19\",\"24+39=63\",\"60-53=7\",\"26+69=95\",\"49+25=74\",\"90-44=46\",\"92-23=69\",\"9+75=84\",\"9+74=83\",\"35-26=9\",\"90-35=55\",\"84-59=25\",\"94-16=78\",\"47-9=38\",\"28+68=96\",\"78+9=87\",\"76-29=47\",\"51-18=33\",\"27+46=73\",\"38+37=75\",\"57+36=93\",\"2+89=91\",\"57+7=64\",\"93-39=54\",\"60-5=55\",\"66+25=91\",\"43-4=39\",\"26+59=85\",\"25+17=42\",\"90-41=49\",\"47-18=29\",\"53-49=4\",\"37+57=94\",\"58+24=82\",\"26+69=95\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$numRows = $t.Rows.Count\n$numCols = $t.Columns.Count\n\n# Assigning Cell(r,c).Range.Text replaces only the cell's content run,\n# leaving the run/paragraph formatting (font, size) and the cell-end\n# markers untouched.\n$idx = 0\nfor ($r = 1; $r -le $numRows; $r++) {\n    for ($c = 1; $c -le $numCols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
